# Update "想去人数" (interested-count) figures for the latest generated
# gh-pages output. Two sheets share the same event list data (展览 and
# 全部类型); both need their F-column counts bumped.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (rows are offset by one vs. "全部类型" because it lacks the
# extra leading event row present in "全部类型").
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value  = 321
$wsExhibit.Range("F6").Value  = 1037
$wsExhibit.Range("F7").Value  = 10758
$wsExhibit.Range("F11").Value = 1041
$wsExhibit.Range("F13").Value = 12078
$wsExhibit.Range("F14").Value = 12539

# Sheet "全部类型" (all event types combined).
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value  = 321
$wsAll.Range("F7").Value  = 1037
$wsAll.Range("F8").Value  = 10758
$wsAll.Range("F12").Value = 1041
$wsAll.Range("F14").Value = 12078
$wsAll.Range("F15").Value = 12539

$wb.Save()
